$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update route_ids text values: "purchase." -> "purchase_stock." prefix
#    (covers both the lone reference and the comma-joined
#    "purchase.route_warehouse0_buy,stock.route_warehouse0_mto" reference,
#    since the replace matches the shared substring in both)
$null = $ws.Cells.Replace("purchase.route_warehouse0_buy", "purchase_stock.route_warehouse0_buy")

# 2. Bump row 29 height slightly (12.8 -> 13.8)
$ws.Rows.Item(29).RowHeight = 13.8

# 3. Materialize an explicit, empty, default-styled cell at Q29
#    (it previously had no cell entry at all for that row)
$ws.Range("Q29").Font.Name = "Calibri"
$ws.Range("Q29").Font.Size = 11

# 4. Update the view/selection state of the sheet: scroll the frozen-pane
#    view over towards columns K/Q and leave column Q selected
$null = $ws.Range("A1").Select()
$null = $ws.Range("K1").Select()
$null = $ws.Range("A2").Select()
$null = $ws.Range("Q:Q").Select()
